$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the paragraphs we need to work with by their known text, so the
# script is resilient to any paragraph renumbering.
# ---------------------------------------------------------------------------
function Get-ParaIndexByText($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.StartsWith($needle)) {
            return $i
        }
    }
    return -1
}

$adibIdx  = Get-ParaIndexByText '{"username":"adib"'
$adiblIdx = Get-ParaIndexByText '{"username":"adibl"'

if ($adibIdx -eq -1 -or $adiblIdx -eq -1) {
    throw "Could not locate the sample record paragraphs to edit."
}

# ---------------------------------------------------------------------------
# 1) Remove the (hidden) _GoBack bookmark from its current location -- it
#    currently sits right before the "_id" paragraph's text and needs to
#    move to the new last record paragraph created below.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2) Drop the old "adibl" sample-record paragraph entirely; its slot is
#    being replaced by the two templated record paragraphs inserted below.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item($adiblIdx).Range.Delete()

# ---------------------------------------------------------------------------
# 3) Replace the old "adib" sample-record paragraph with two new record
#    paragraphs: a concrete example record (username1/password1/program1)
#    and a generic templated record (usernameN/passwordN/programN), the
#    latter carrying the restored _GoBack bookmark right before its
#    opening brace.
# ---------------------------------------------------------------------------
$xmlWrapperOpen = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
'@

$xmlWrapperClose = @'
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$recordsXml = $xmlWrapperOpen + @'
<w:p>
  <w:r><w:tab/></w:r>
  <w:r><w:t>{"username":"</w:t></w:r>
  <w:r><w:t>username1</w:t></w:r>
  <w:r><w:t>","password":"</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t>password</w:t></w:r>
  <w:r><w:t>1</w:t></w:r>
  <w:r><w:t>","program_id":"</w:t></w:r>
  <w:r><w:t>program1</w:t></w:r>
  <w:r><w:t>"},</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:tab/></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r><w:t>{"username":"</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>username</w:t></w:r>
  <w:r><w:t>N</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>","password":"</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>password</w:t></w:r>
  <w:r><w:t>N</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>","</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>program_id</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>":"</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>program</w:t></w:r>
  <w:r><w:t>N</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>"}</w:t></w:r>
</w:p>
'@ + $xmlWrapperClose

$adibRange = $d.Paragraphs.Item($adibIdx).Range
[void]$adibRange.InsertXML($recordsXml)
